# Edit: normalize Polish-diacritic vehicle-segment / range labels on
# "Arkusz5" to their plain-ASCII equivalents, and leave that sheet as the
# active sheet/selection (as it was when the workbook was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz5")

# Header: "zasięg" -> "zasieg"
$ws.Range("B1").Value = "zasieg"

# Column A body values: strip Polish diacritics where the author retyped them
$ws.Range("A2").Value  = "male"
$ws.Range("A3").Value  = "kombi-van(male)"
$ws.Range("A5").Value  = "kombi-van(srednie)"
$ws.Range("A7").Value  = "dostawcze(male)"
$ws.Range("A9").Value  = "dostawcze(srednie)"
$ws.Range("A10").Value = "klasa wyzsza-srednia"
$ws.Range("A11").Value = "dostawcze(duze)"

# Leave Arkusz5 as the active sheet with A11 selected, matching the
# workbook state at the time of last save.
$ws.Activate()
$ws.Range("A11").Select()
